# Refactor cost income rule
# - Rename the cost/income "Project Y" / "Project Z" rules to "Project A" / "Project B"
# - Add two new cost/income rules: "Project C" and "Project D"
# - Update the underlying dates for the renamed rules
# - Leave the selection/active-sheet state matching the author's final position

$wb = $excel.ActiveWorkbook

# --- other sheets: restore each sheet's own last-used selection -----------------
# (done first, since selecting a range on a sheet also activates that sheet;
#  "costs" must be the very last sheet activated below)
$production = $wb.Worksheets.Item("production")
[void]$production.Range("A9").Select()

$dividend = $wb.Worksheets.Item("dividend")
[void]$dividend.Range("C13").Select()

# --- "costs" sheet (CostIncome) -------------------------------------------------
$costs = $wb.Worksheets.Item("costs")

# Row 4: was "Project Y", now "Project A" with refreshed dates
$costs.Range("A4").Value = "2025-01-01"
$costs.Range("B4").Value = "Project A"
$costs.Range("C4").Value = -100
$costs.Range("D4").Value = "2025-01-15"
$costs.Range("E4").Value = "2025-09-30"

# Row 5: was "Project Z", now "Project B" with refreshed dates
$costs.Range("A5").Value = "2025-11-30"
$costs.Range("B5").Value = "Project B"
$costs.Range("C5").Value = 10000
$costs.Range("D5").Value = "2026-02-15"
$costs.Range("E5").Value = "2026-05-16"

# Row 6 (new): "Project C"
$costs.Range("A6").Value = "2026-05-30"
$costs.Range("B6").Value = "Project C"
$costs.Range("C6").Value = -500
$costs.Range("D6").Value = "2025-01-01"
$costs.Range("E6").Value = "2025-12-31"

# Row 7 (new): "Project D"
$costs.Range("A7").Value = "2026-03-31"
$costs.Range("B7").Value = "Project D"
$costs.Range("C7").Value = 600
$costs.Range("D7").Value = "2025-01-01"
$costs.Range("E7").Value = "2025-12-31"

# Make "costs" the active sheet/tab and leave the cursor on the last new row,
# matching the saved view state of the edited workbook.
$costs.Activate()
[void]$costs.Range("A7").Select()
